# Auto-generated edit script: update cryptos list prices/volumes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.893.56"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.813.66"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.41"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4290"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3691"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07243"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8621"
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.21"
$ws.Range("E11").Value = "  +4.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.012.85"
$ws.Range("E12").Value = "  +11.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.641"
$ws.Range("E13").Value = "  +4.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.397"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06896"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "80.73"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008937"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.18"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.942.91"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.182"
$ws.Range("E22").Value = "  +2.73%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.234.87"
$ws.Range("E24").Value = "  +10.13%  "
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.28"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.226"
$ws.Range("E28").Value = "  +4.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.96"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  +14.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08941"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7441"
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.167"
$ws.Range("E33").Value = "  +7.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.424"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.798"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.124"
$ws.Range("E37").Value = "  +3.97%  "
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5078"
$ws.Range("E40").Value = "  +2.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1640"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("E42").Value = "  +8.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.436"
$ws.Range("E43").Value = "  +7.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.253"
$ws.Range("E44").Value = "  +4.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.74"
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.40"
$ws.Range("E46").Value = "  +3.50%  "
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.655"
$ws.Range("E48").Value = "  +5.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06287"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4553"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.806"
$ws.Range("E51").Value = "  +5.96%  "
